# Fixing missing citys: Rio Grande and Loiza
# Update the 2020 COVID case counts in column C now that Rio Grande and
# Loiza (both part of the Fajardo region) are included in the totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4100
$ws.Range("C3").Value = 2114
$ws.Range("C4").Value = 2070
$ws.Range("C5").Value = 1883
$ws.Range("C6").Value = 1243

$ws.Range("C9").Value = 7817
$ws.Range("C10").Value = 4177
$ws.Range("C11").Value = 4050
$ws.Range("C12").Value = 3836
$ws.Range("C13").Value = 2296

$ws.Range("C15").Value = 638
$ws.Range("C16").Value = 5480
$ws.Range("C17").Value = 2804
$ws.Range("C18").Value = 2858
$ws.Range("C19").Value = 2401
$ws.Range("C20").Value = 1589
$ws.Range("C21").Value = 887
$ws.Range("C22").Value = 476

# Fajardo region (includes Rio Grande and Loiza) - previously missing data
$ws.Range("C23").Value = 761
$ws.Range("C24").Value = 414
$ws.Range("C25").Value = 431
$ws.Range("C26").Value = 433
$ws.Range("C27").Value = 240
$ws.Range("C28").Value = 131
$ws.Range("C29").Value = 46

$ws.Range("C30").Value = 3160
$ws.Range("C32").Value = 1742
$ws.Range("C33").Value = 1526
$ws.Range("C34").Value = 930

$ws.Range("C37").Value = 9878
$ws.Range("C38").Value = 5674
$ws.Range("C39").Value = 5480
$ws.Range("C40").Value = 5311
$ws.Range("C41").Value = 3092
$ws.Range("C42").Value = 1727
$ws.Range("C43").Value = 1037

$ws.Range("C44").Value = 2969
$ws.Range("C45").Value = 1501
$ws.Range("C46").Value = 1632
$ws.Range("C47").Value = 1365
$ws.Range("C48").Value = 822

# Restore the view state saved with the workbook: zoomed in, scrolled down
# to the Fajardo rows, with C49 as the active selected cell.
$ws.Range("A32").Select()
$excel.ActiveWindow.Zoom = 158
$ws.Range("C49").Select()
